function Set-TextValue($ws, $ref, $val) {
    # Preserve the cell's existing style while forcing the written value to
    # be stored as text (not auto-converted to a number/date by Excel).
    $rng = $ws.Range($ref)
    $orig = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $orig
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "66.652.78"
Set-TextValue $ws "E2" "  -5.94%  "
Set-TextValue $ws "D3" "3.208.70"
Set-TextValue $ws "E3" "  -9.18%  "
Set-TextValue $ws "E4" "  -0.04%  "
Set-TextValue $ws "D5" "576.98"
Set-TextValue $ws "E5" "  -5.97%  "
Set-TextValue $ws "D6" "150.50"
Set-TextValue $ws "E6" "  -13.53%  "
Set-TextValue $ws "D7" "0.999"
Set-TextValue $ws "E7" "  -0.10%  "
Set-TextValue $ws "D8" "3.199.88"
Set-TextValue $ws "E8" "  -9.33%  "
Set-TextValue $ws "D9" "0.540"
Set-TextValue $ws "E9" "  -11.55%  "
Set-TextValue $ws "E10" "  -13.29%  "
Set-TextValue $ws "D11" "6.80"
Set-TextValue $ws "E11" "  -8.23%  "
Set-TextValue $ws "D12" "0.497"
Set-TextValue $ws "E12" "  -15.70%  "
Set-TextValue $ws "D13" "38.19"
Set-TextValue $ws "E14" "  -12.50%  "
Set-TextValue $ws "D15" "3.724.46"
Set-TextValue $ws "E15" "  -9.13%  "
Set-TextValue $ws "D16" "66.615.29"
Set-TextValue $ws "E16" "  -5.92%  "
Set-TextValue $ws "B17" "BitcoinCash"
Set-TextValue $ws "C17" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws "D17" "541.79"
Set-TextValue $ws "E17" "  -12.10%  "
Set-TextValue $ws "B18" "WrappedEther"
Set-TextValue $ws "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D18" "3.208.03"
Set-TextValue $ws "E18" "  -9.05%  "
Set-TextValue $ws "D20" "7.09"
Set-TextValue $ws "E20" "  -16.05%  "
Set-TextValue $ws "D21" "15.00"
Set-TextValue $ws "E21" "  -15.61%  "
Set-TextValue $ws "D22" "0.752"
Set-TextValue $ws "E22" "  -15.18%  "
Set-TextValue $ws "D23" "7.69"
Set-TextValue $ws "E23" "  -14.61%  "
Set-TextValue $ws "D24" "85.13"
Set-TextValue $ws "E24" "  -13.33%  "
Set-TextValue $ws "D25" "13.32"
Set-TextValue $ws "E25" "  -15.57%  "
Set-TextValue $ws "E26" "  -0.02%  "
Set-TextValue $ws "D27" "3.12"
Set-TextValue $ws "E27" "  -17.62%  "
Set-TextValue $ws "D28" "8.02"
Set-TextValue $ws "E28" "  -12.40%  "
Set-TextValue $ws "D29" "29.04"
Set-TextValue $ws "E29" "  -14.42%  "
Set-TextValue $ws "E30" "  -19.14%  "
Set-TextValue $ws "D31" "2.55"
Set-TextValue $ws "E31" "  -16.11%  "
Set-TextValue $ws "E32" "  -13.37%  "
Set-TextValue $ws "D33" "536.14"
Set-TextValue $ws "E33" "  -12.96%  "
Set-TextValue $ws "D34" "6.49"
Set-TextValue $ws "E34" "  -20.65%  "
Set-TextValue $ws "E35" "  -17.66%  "
Set-TextValue $ws "E36" "  +0.06%  "
Set-TextValue $ws "D37" "52.85"
Set-TextValue $ws "E37" "  -7.44%  "
Set-TextValue $ws "D38" "0.0422"
Set-TextValue $ws "E38" "  -11.21%  "
Set-TextValue $ws "D39" "0.0840"
Set-TextValue $ws "E39" "  -16.51%  "
Set-TextValue $ws "D40" "9.08"
Set-TextValue $ws "E40" "  -16.41%  "
Set-TextValue $ws "D41" "0.123"
Set-TextValue $ws "E41" "  -14.52%  "
Set-TextValue $ws "D42" "2.911.03"
Set-TextValue $ws "E42" "  -13.92%  "
Set-TextValue $ws "E43" "  -27.38%  "
Set-TextValue $ws "B44" "TheGraph"
Set-TextValue $ws "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D44" "0.260"
Set-TextValue $ws "E44" "  -17.40%  "
Set-TextValue $ws "B45" "PEPE"
Set-TextValue $ws "C45" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws "D45" "0.0₃0582"
Set-TextValue $ws "E45" "  -21.60%  "
Set-TextValue $ws "D47" "2.34"
Set-TextValue $ws "E47" "  -20.81%  "
Set-TextValue $ws "D48" "25.78"
Set-TextValue $ws "E48" "  -20.34%  "
Set-TextValue $ws "E49" "  -19.11%  "
Set-TextValue $ws "B50" "Stellar"
Set-TextValue $ws "C50" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D50" "0.112"
Set-TextValue $ws "E50" "  -13.96%  "
Set-TextValue $ws "B51" "Monero"
Set-TextValue $ws "C51" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D51" "123.17"
Set-TextValue $ws "E51" "  -8.16%  "
